# Rebuild the "Hoja1" sheet with the new rut/nombre/estado/prevision/email/afiliado table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- clear the old demo data (A1:B5) so no stray cells are left behind ---
$ws.Range("A1:B5").ClearContents()

# --- header row ---
$ws.Range("A1").Value = "rut"
$ws.Range("B1").Value = "nombre"
$ws.Range("C1").Value = "estado"
$ws.Range("D1").Value = "prevision"
$ws.Range("E1").Value = "email"
$ws.Range("F1").Value = "afiliado"

# --- data rows ---
# The "afiliado" column holds literal text "true"/"false", not Excel booleans.
# Typing the bare word auto-coerces it to a real Boolean, so enter it with a
# leading apostrophe to force text, then strip the resulting quote-prefix
# style so the cell ends up as a plain, unstyled shared-string cell.

# row 2
$ws.Range("A2").Value = 82937288
$ws.Range("B2").Value = "yerso"
$ws.Range("C2").Value = "activo"
$ws.Range("D2").Value = "fonasa"
$ws.Range("E2").Value = "hola"
$ws.Range("F2").Value = "'true"

# row 3
$ws.Range("A3").Value = "8237478-9"
$ws.Range("B3").Value = "thomas"
$ws.Range("C3").Value = "inactivo"
$ws.Range("D3").Value = "isapre"
$ws.Range("E3").Value = "culi"
$ws.Range("F3").Value = "'false"

# row 4
$ws.Range("A4").Value = 1248238
$ws.Range("B4").Value = "leandro"
$ws.Range("C4").Value = "activo"
$ws.Range("D4").Value = "isapre"
$ws.Range("E4").Value = "ctm"
$ws.Range("F4").Value = "'true"

# row 5
$ws.Range("A5").Value = 23743728
$ws.Range("B5").Value = "perkinazo"
$ws.Range("C5").Value = "inactivo"
$ws.Range("D5").Value = "fonasa"
$ws.Range("E5").Value = "asd"
$ws.Range("F5").Value = "'false"

# row 6
$ws.Range("A6").Value = 4357984
$ws.Range("B6").Value = "jorsi"
$ws.Range("C6").Value = "inactivo"
$ws.Range("D6").Value = "fonasa"
$ws.Range("E6").Value = "yapo"
$ws.Range("F6").Value = "'true"

# row 7
$ws.Range("A7").Value = 12345678
$ws.Range("B7").Value = "nelson"
$ws.Range("C7").Value = "activo"
$ws.Range("D7").Value = "genial"
$ws.Range("E7").Value = "bueno"
$ws.Range("F7").Value = "'false"

# Strip the quote-prefix style the apostrophe entry above added.
$ws.Range("F2:F7").Style = "Normal"

# --- page setup: A4 paper, portrait orientation ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- selection ends on the last cell entered ---
[void]$ws.Range("F7").Select()
